# "changes to the Create Datesets for Assertion"
# - Flip the Exceute flag for the "Create Datasets" test case (row 3) from N -> Y
# - Flip the Exceute flag for the "Create Input Forms" test case (row 7) from Y -> N
# - Move the active selection on the Test Cases sheet from D8 to D5

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

$ws.Range("D3").Value = "Y"
$ws.Range("D7").Value = "N"

$ws.Activate()
$ws.Range("D5").Select()
